# Apropriações Mauro.xlsx — add the next time-tracking entry (row 3):
#   A3 = 07/10/2013 (date, same format as A2)
#   B3 = 4:50 (duration, same format as B2)
# and extend the selection to cover B2:B3 (matching the post-edit sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
$ws.Range("A3").Value = 41554
$ws.Range("B3").Value = 0.20138888888888887

# Copy the formatting (number format / style) from row 2 down to row 3 so
# the new cells reuse the existing date/time cell styles instead of
# creating new ones.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats

# --- Selection ----------------------------------------------------------
# Extend the current selection down to include the newly added row.
$ws.Range("B2:B3").Select()
